$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.957.64"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "2.018.40"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'226.62"
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").Value = "'0.608"
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'54.97"
$ws.Range("E8").Value = "  -4.95%  "
$ws.Range("D9").Value = "'0.380"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").Value = "'0.0792"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  -3.72%  "
$ws.Range("D12").Value = "2.316.65"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "'14.29"
$ws.Range("E13").Value = "  -5.08%  "
$ws.Range("D14").Value = "'20.58"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").Value = "'0.745"
$ws.Range("E15").Value = "  -3.38%  "
$ws.Range("D16").Value = "'5.15"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("D17").Value = "2.039.36"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "36.872.73"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "'6.09"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "'68.80"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").Value = "'226.65"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("D26").Value = "'166.99"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("E27").Value = "  -3.96%  "
$ws.Range("E28").Value = "  -4.72%  "
$ws.Range("D29").Value = "'18.75"
$ws.Range("E29").Value = "  -3.86%  "
$ws.Range("E30").Value = "  -2.96%  "
$ws.Range("E31").Value = "  -4.41%  "
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").Value = "'4.45"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'3.18"
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("D39").Value = "'5.39"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").Value = "'0.0219"
$ws.Range("E40").Value = "  -5.55%  "
$ws.Range("D41").Value = "1.490.26"
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("D42").Value = "'17.07"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").Value = "'0.0930"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("D44").Value = "'95.00"
$ws.Range("E44").Value = "  -5.32%  "
$ws.Range("D45").Value = "'2.80"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("E46").Value = "  -5.21%  "
$ws.Range("D47").Value = "'7.28"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("D49").Value = "'2.90"
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("E50").Value = "  -5.99%  "
$ws.Range("D51").Value = "2.207.02"
$ws.Range("E51").Value = "  -2.90%  "
